$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A9 already carries style index 4 ("absoluteStyle"); copy just the
# formatting (not the value/content) from A9 onto A10, A11 and A12 so they
# pick up the same look without altering their existing text.
$ws.Range("A9").Copy()
$ws.Range("A10").PasteSpecial(-4122)
$ws.Range("A11").PasteSpecial(-4122)
$ws.Range("A12").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# C11 must keep holding its number-as-text value ("-1" instead of "-3").
# Force text storage, then restore C11's original look (border/font/
# alignment) by pasting format from the sibling cell C10, which still
# carries the untouched style.
$ws.Range("C11").NumberFormat = "@"
$ws.Range("C11").Value = "-1"
$ws.Range("C10").Copy()
$ws.Range("C11").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Plain numeric + textual updates.
$ws.Range("C12").Value = -11
$ws.Range("E12").Value = "24/140"
